$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-parsed as a number;
# force text storage (matches the inlineStr cells in the source workbook).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated Price (D) / Volume(1h) (E) values row by row.
$ws.Range("D2").Value = '64.382.64'
$ws.Range("E2").Value = '  +1.44%  '
$ws.Range("D3").Value = '2.620.13'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '593.55'
$ws.Range("E5").Value = '  -0.27%  '
$ws.Range("D6").Value = '151.80'
$ws.Range("E6").Value = '  +1.07%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  +5.25%  '
$ws.Range("D10").Value = '0.397'
$ws.Range("D11").Value = '5.80'
$ws.Range("E11").Value = '  +2.11%  '
$ws.Range("E12").Value = '  +1.34%  '
$ws.Range("D13").Value = '28.52'
$ws.Range("E13").Value = '  +3.14%  '
$ws.Range("D14").Value = '3.087.20'
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("D15").Value = '0.0000172'
$ws.Range("E15").Value = '  +13.15%  '
$ws.Range("D16").Value = '64.274.71'
$ws.Range("E16").Value = '  +1.43%  '
$ws.Range("D17").Value = '2.561.40'
$ws.Range("E17").Value = '  -1.83%  '
$ws.Range("E18").Value = '  -0.34%  '
$ws.Range("D19").Value = '4.79'
$ws.Range("E19").Value = '  +2.41%  '
$ws.Range("D20").Value = '350.02'
$ws.Range("E20").Value = '  +0.98%  '
$ws.Range("D21").Value = '7.15'
$ws.Range("E21").Value = '  +4.45%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.29%  '
$ws.Range("D23").Value = '67.57'
$ws.Range("E23").Value = '  +1.47%  '
$ws.Range("E24").Value = '  -0.89%  '
$ws.Range("D25").Value = '9.30'
$ws.Range("E25").Value = '  +1.25%  '
$ws.Range("D26").Value = '1.65'
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("E27").Value = '  +1.60%  '
$ws.Range("E28").Value = '  +1.58%  '
$ws.Range("D29").Value = '544.44'
$ws.Range("E29").Value = '  -1.80%  '
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").Value = '0.0₃0909'
$ws.Range("E31").Value = '  +7.63%  '
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("D33").Value = '1.82'
$ws.Range("E33").Value = '  +4.17%  '
$ws.Range("D34").Value = '5.66'
$ws.Range("E34").Value = '  +8.09%  '
$ws.Range("D35").Value = '6.23'
$ws.Range("E35").Value = '  +1.23%  '
$ws.Range("E36").Value = '  +2.58%  '
$ws.Range("D37").Value = '163.79'
$ws.Range("E37").Value = '  -2.33%  '
$ws.Range("D38").Value = '20.12'
$ws.Range("E38").Value = '  +3.35%  '
$ws.Range("D39").Value = '2.00'
$ws.Range("E39").Value = '  +3.73%  '
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").Value = '168.36'
$ws.Range("E42").Value = '  +1.05%  '
$ws.Range("D43").Value = '41.63'
$ws.Range("E43").Value = '  +4.91%  '
$ws.Range("D44").Value = '4.09'
$ws.Range("E44").Value = '  +4.88%  '
$ws.Range("D45").Value = '23.22'
$ws.Range("E45").Value = '  +7.52%  '
$ws.Range("D46").Value = '0.0597'
$ws.Range("E46").Value = '  +2.26%  '
$ws.Range("D47").Value = '2.20'
$ws.Range("E47").Value = '  +10.65%  '
$ws.Range("E48").Value = '  +1.76%  '
$ws.Range("D49").Value = '0.0251'
$ws.Range("E49").Value = '  +0.41%  '
$ws.Range("D50").Value = '0.0980'
$ws.Range("E50").Value = '  +1.60%  '
$ws.Range("D51").Value = '19.33'
$ws.Range("E51").Value = '  +0.86%  '

# Restore the default (unstyled) formatting now that the text is locked in.
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
